$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3g"
$ws.Range("C2").Value = "Nrp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.58860066666667
$ws.Range("H2").Value = 40.765802
$ws.Range("I2").Value = 0.9060457790710231
$ws.Range("J2").Value = 0.9060457790710231
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 70.23436
$ws.Range("N2").Value = 210.70308
$ws.Range("O2").Value = 0.7023186840741513
$ws.Range("P2").Value = 0.7023186840741513
$ws.Range("Q2").Value = 954.3866711189066
$ws.Range("R2").Value = 8589.48004007016
$ws.Range("S2").Value = 0.6363328792681002
$ws.Range("T2").Value = 0.6363328792681002

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3g"
$ws.Range("C3").Value = "Nrp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.58860066666667
$ws.Range("H3").Value = 40.765802
$ws.Range("I3").Value = 0.9060457790710231
$ws.Range("J3").Value = 0.9060457790710231
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.551362
$ws.Range("N3").Value = 19.654086
$ws.Range("O3").Value = 0.06551129587759326
$ws.Range("P3").Value = 0.06551129587759325
$ws.Range("Q3").Value = 89.02384204077467
$ws.Range("R3").Value = 801.214578366972
$ws.Range("S3").Value = 0.05935623311136629
$ws.Range("T3").Value = 0.05935623311136628

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3g"
$ws.Range("C4").Value = "Nrp2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 13.58860066666667
$ws.Range("H4").Value = 40.765802
$ws.Range("I4").Value = 0.9060457790710231
$ws.Range("J4").Value = 0.9060457790710231
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.268944666666666
$ws.Range("N4").Value = 12.806834
$ws.Range("O4").Value = 0.04268793224112385
$ws.Range("P4").Value = 0.04268793224112385
$ws.Range("Q4").Value = 58.00898434342977
$ws.Range("R4").Value = 522.080859090868
$ws.Range("S4").Value = 0.03867722082434011
$ws.Range("T4").Value = 0.03867722082434011

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema3g"
$ws.Range("C5").Value = "Nrp2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 13.58860066666667
$ws.Range("H5").Value = 40.765802
$ws.Range("I5").Value = 0.9060457790710231
$ws.Range("J5").Value = 0.9060457790710231
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 18.948881
$ws.Range("N5").Value = 56.846643
$ws.Range("O5").Value = 0.1894820878071316
$ws.Range("P5").Value = 0.1894820878071315
$ws.Range("Q5").Value = 257.4887769891873
$ws.Range("R5").Value = 2317.398992902686
$ws.Range("S5").Value = 0.1716794458672165
$ws.Range("T5").Value = 0.1716794458672165

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3g"
$ws.Range("C6").Value = "Nrp2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05898900000000001
$ws.Range("H6").Value = 0.176967
$ws.Range("I6").Value = 0.003933203703066158
$ws.Range("J6").Value = 0.003933203703066157
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 70.23436
$ws.Range("N6").Value = 210.70308
$ws.Range("O6").Value = 0.7023186840741513
$ws.Range("P6").Value = 0.7023186840741513
$ws.Range("Q6").Value = 4.14305466204
$ws.Range("R6").Value = 37.28749195836
$ws.Range("S6").Value = 0.002762362448933003
$ws.Range("T6").Value = 0.002762362448933002

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3g"
$ws.Range("C7").Value = "Nrp2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05898900000000001
$ws.Range("H7").Value = 0.176967
$ws.Range("I7").Value = 0.003933203703066158
$ws.Range("J7").Value = 0.003933203703066157
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.551362
$ws.Range("N7").Value = 19.654086
$ws.Range("O7").Value = 0.06551129587759326
$ws.Range("P7").Value = 0.06551129587759325
$ws.Range("Q7").Value = 0.3864582930180001
$ws.Range("R7").Value = 3.478124637162
$ws.Range("S7").Value = 0.0002576692715384126
$ws.Range("T7").Value = 0.0002576692715384124

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema3g"
$ws.Range("C8").Value = "Nrp2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.05898900000000001
$ws.Range("H8").Value = 0.176967
$ws.Range("I8").Value = 0.003933203703066158
$ws.Range("J8").Value = 0.003933203703066157
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.268944666666666
$ws.Range("N8").Value = 12.806834
$ws.Range("O8").Value = 0.04268793224112385
$ws.Range("P8").Value = 0.04268793224112385
$ws.Range("Q8").Value = 0.251820776942
$ws.Range("R8").Value = 2.266386992478
$ws.Range("S8").Value = 0.0001679003331670256
$ws.Range("T8").Value = 0.0001679003331670255

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema3g"
$ws.Range("C9").Value = "Nrp2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.05898900000000001
$ws.Range("H9").Value = 0.176967
$ws.Range("I9").Value = 0.003933203703066158
$ws.Range("J9").Value = 0.003933203703066157
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 18.948881
$ws.Range("N9").Value = 56.846643
$ws.Range("O9").Value = 0.1894820878071316
$ws.Range("P9").Value = 0.1894820878071315
$ws.Range("Q9").Value = 1.117775541309
$ws.Range("R9").Value = 10.059979871781
$ws.Range("S9").Value = 0.0007452716494277168
$ws.Range("T9").Value = 0.0007452716494277165

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Sema3g"
$ws.Range("C10").Value = "Nrp2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.350108
$ws.Range("H10").Value = 4.050324
$ws.Range("I10").Value = 0.09002101722591065
$ws.Range("J10").Value = 0.09002101722591065
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 70.23436
$ws.Range("N10").Value = 210.70308
$ws.Range("O10").Value = 0.7023186840741513
$ws.Range("P10").Value = 0.7023186840741513
$ws.Range("Q10").Value = 94.82397131087998
$ws.Range("R10").Value = 853.41574179792
$ws.Range("S10").Value = 0.06322344235711808
$ws.Range("T10").Value = 0.06322344235711808

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Sema3g"
$ws.Range("C11").Value = "Nrp2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.350108
$ws.Range("H11").Value = 4.050324
$ws.Range("I11").Value = 0.09002101722591065
$ws.Range("J11").Value = 0.09002101722591065
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.551362
$ws.Range("N11").Value = 19.654086
$ws.Range("O11").Value = 0.06551129587759326
$ws.Range("P11").Value = 0.06551129587759325
$ws.Range("Q11").Value = 8.845046247095999
$ws.Range("R11").Value = 79.605416223864
$ws.Range("S11").Value = 0.005897393494688552
$ws.Range("T11").Value = 0.00589739349468855

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Sema3g"
$ws.Range("C12").Value = "Nrp2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.350108
$ws.Range("H12").Value = 4.050324
$ws.Range("I12").Value = 0.09002101722591065
$ws.Range("J12").Value = 0.09002101722591065
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.268944666666666
$ws.Range("N12").Value = 12.806834
$ws.Range("O12").Value = 0.04268793224112385
$ws.Range("P12").Value = 0.04268793224112385
$ws.Range("Q12").Value = 5.763536346023999
$ws.Range("R12").Value = 51.87182711421599
$ws.Range("S12").Value = 0.003842811083616717
$ws.Range("T12").Value = 0.003842811083616717

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Sema3g"
$ws.Range("C13").Value = "Nrp2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.350108
$ws.Range("H13").Value = 4.050324
$ws.Range("I13").Value = 0.09002101722591065
$ws.Range("J13").Value = 0.09002101722591065
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 18.948881
$ws.Range("N13").Value = 56.846643
$ws.Range("O13").Value = 0.1894820878071316
$ws.Range("P13").Value = 0.1894820878071315
$ws.Range("Q13").Value = 25.583035829148
$ws.Range("R13").Value = 230.247322462332
$ws.Range("S13").Value = 0.01705737029048731
$ws.Range("T13").Value = 0.0170573702904873

# Remove the old Resolving-Mac sending-cluster rows (14-17)
$ws.Range("A14:T17").Delete()
